$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:K1 (copy formatting from an existing header cell, then set text) ---
$ws.Range("A1").Copy()
$ws.Range("F1:K1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "C4.5 acc"
$ws.Range("G1").Value = "credal-C4.5 acc"
$ws.Range("H1").Value = "SPN acc"
$ws.Range("I1").Value = "CSPN low"
$ws.Range("J1").Value = "CSPN high"
$ws.Range("K1").Value = "CSPN robust"

# --- Row 2 (A2 = 0) ---
$ws.Range("B2").Value = 69.88888888888889
$ws.Range("C2").Value = 69.44444444444444
$ws.Range("D2").Value = 70.31111111111112
$ws.Range("E2").Value = 70.05007533305347
$ws.Range("F2").Value = 64.82222222222222
$ws.Range("G2").Value = 65.03333333333333
$ws.Range("H2").Value = 69.67777777777778
$ws.Range("I2").Value = 69.67777777777778
$ws.Range("J2").Value = 69.67777777777778
$ws.Range("K2").Value = 69.67777777777778

# --- Row 3 (A3 = 5) ---
$ws.Range("B3").Value = 70.52222222222223
$ws.Range("C3").Value = 70.12222222222222
$ws.Range("D3").Value = 70.84444444444445
$ws.Range("E3").Value = 70.63313254428461
$ws.Range("F3").Value = 58.01111111111111
$ws.Range("G3").Value = 59.68888888888888
$ws.Range("H3").Value = 70.04444444444445
$ws.Range("I3").Value = 56.33333333333334
$ws.Range("J3").Value = 79.14444444444445
$ws.Range("K3").Value = 73.18615626658946

# --- Row 4 (A4 = 10) ---
$ws.Range("B4").Value = 70.83333333333333
$ws.Range("C4").Value = 70.31111111111112
$ws.Range("D4").Value = 71.44444444444444
$ws.Range("E4").Value = 71.11849361731574
$ws.Range("F4").Value = 51.17777777777778
$ws.Range("G4").Value = 58.37777777777777
$ws.Range("H4").Value = 70.47777777777779
$ws.Range("I4").Value = 37.47777777777777
$ws.Range("J4").Value = 89.16666666666667
$ws.Range("K4").Value = 77.9314636777633

# --- Row 5 (A5 = 20) ---
$ws.Range("B5").Value = 70.14444444444445
$ws.Range("C5").Value = 69.46666666666667
$ws.Range("D5").Value = 70.78888888888888
$ws.Range("E5").Value = 70.41087747269458
$ws.Range("F5").Value = 37.2
$ws.Range("G5").Value = 44.63333333333333
$ws.Range("H5").Value = 72.2
$ws.Range("I5").Value = 25.28888888888889
$ws.Range("J5").Value = 91.75555555555555
$ws.Range("K5").Value = 74.90104613419182

# --- Row 6 (A6 = 30) ---
$ws.Range("B6").Value = 70.07777777777778
$ws.Range("C6").Value = 69.37777777777778
$ws.Range("D6").Value = 70.8
$ws.Range("E6").Value = 70.38341182815928
$ws.Range("F6").Value = 35.67777777777778
$ws.Range("G6").Value = 44.36666666666667
$ws.Range("H6").Value = 72.90000000000001
$ws.Range("I6").Value = 28.93333333333333
$ws.Range("J6").Value = 88.88888888888889
$ws.Range("K6").Value = 71.57912578987536
